# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.135.30"
$ws.Range("E2").Value = "  -2.20%  "
$ws.Range("D3").Value = "1.852.18"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'237.79"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "'0.6889"
$ws.Range("E6").Value = "  -5.46%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.07774"
$ws.Range("E8").Value = "  +8.96%  "
$ws.Range("D9").Value = "'0.3036"
$ws.Range("E9").Value = "  -3.09%  "
$ws.Range("D10").Value = "'23.23"
$ws.Range("E10").Value = "  -4.94%  "
$ws.Range("D11").Value = "'0.08138"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "1.852.47"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").Value = "'0.7247"
$ws.Range("E13").Value = "  -2.43%  "
$ws.Range("D14").Value = "'5.203"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "'89.04"
$ws.Range("E15").Value = "  -3.64%  "
$ws.Range("D16").Value = "29.136.85"
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.000007825"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'5.737"
$ws.Range("E18").Value = "  -4.30%  "
$ws.Range("D19").Value = "'13.19"
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("D20").Value = "'235.74"
$ws.Range("E20").Value = "  -4.73%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "2.101.17"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").Value = "'0.9998"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'7.588"
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("D25").Value = "'161.62"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").Value = "'8.959"
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("D27").Value = "'0.1421"
$ws.Range("E27").Value = "  -7.36%  "
$ws.Range("D28").Value = "'18.06"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("D29").Value = "'1.967"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("D30").Value = "'1.397"
$ws.Range("E30").Value = "  -3.48%  "
$ws.Range("D31").Value = "'4.534"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").Value = "'1.485"
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("D34").Value = "'0.05200"
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("D35").Value = "'1.179"
$ws.Range("E35").Value = "  -4.11%  "
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "'1.028"
$ws.Range("E36").Value = "  +2.91%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.7030"
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("D38").Value = "'2.652"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("D39").Value = "'0.01853"
$ws.Range("E39").Value = "  -4.24%  "
$ws.Range("D40").Value = "'2.672"
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("D41").Value = "'0.9050"
$ws.Range("E41").Value = "  +3.29%  "
$ws.Range("D42").Value = "1.101.02"
$ws.Range("E42").Value = "  +5.72%  "
$ws.Range("D43").Value = "'6.008"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").Value = "'0.4270"
$ws.Range("E44").Value = "  -4.54%  "
$ws.Range("D45").Value = "'70.47"
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("D46").Value = "'0.9998"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "'103.02"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").Value = "'1.758"
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("D49").Value = "1.998.37"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("D50").Value = "'9.151"
$ws.Range("E50").Value = "  -3.89%  "
$ws.Range("D51").Value = "'6.953"
$ws.Range("E51").Value = "  -6.74%  "
